$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 21: Leetcode Question No. 1721, Question "Swapping Nodes in a Linked List"
$ws.Range("A21").Value = 1721
$ws.Range("A21").HorizontalAlignment = -4131

$ws.Range("B21").Value = "Swapping Nodes in a Linked List"

# Replicate the multi-area selection saved in the file (B21 and G9, active cell G9)
$ws.Range("B21").Select()
$ws.Range("G9").Activate()
